$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: date, category, amount
$ws.Range("A3").Value = 45232
$ws.Range("A3").Style = $ws.Range("A2").Style
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat

$ws.Range("B3").Value = "General donations"

$ws.Range("C3").Value = 3000
